# Updated RAD to add MD Central Registration Number to 2 tax forms.
# The "Date" column (B) on the RAD results sheet is refreshed with the
# latest test-run timestamps (4 new shared strings replace the 4 stale
# ones referenced by B2:B5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Fri Oct 06 11:31:08 EDT 2023"
$ws.Range("B3").Value = "Fri Oct 06 11:31:21 EDT 2023"
$ws.Range("B4").Value = "Fri Oct 06 11:31:34 EDT 2023"
$ws.Range("B5").Value = "Fri Oct 06 11:31:46 EDT 2023"
